# The "Prix Spot" sheet tracks one column of hourly prices per calendar
# day. A new day ("23-dec") was inserted into the date sequence right
# before the existing "01-oct." column (at column ET), pushing every
# column from ET onward one position to the right (ET -> EU, ... ,
# FX -> FY) and growing the sheet's used range from A1:FX25 to A1:FY25.
#
# Inserting a whole column reproduces that shift exactly (values, formats
# and the used-range/dimension all move together), then we just need to
# fill in the freshly inserted column: the header label in row 1 and a
# "-" placeholder (matching the sheet's existing "no data" convention)
# in the 24 hourly data rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

$ws.Columns("ET:ET").Insert()

$ws.Range("ET1").Value = "23-dec"
$ws.Range("ET2:ET25").Value = "-"
